# Update parameter-estimate tables (HR / survival model fits) and
# covariance matrices across the model worksheets.

$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -4.00779465950803
$ws.Range("C2").Value = 0.628093159789506
$ws.Range("B3").Value = 0.375222007563247
$ws.Range("C3").Value = 0.244402793342619

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.88150849650292
$ws.Range("C2").Value = 0.698101806909137
$ws.Range("B3").Value = -1.04537780645883
$ws.Range("C3").Value = 0.221833270937624

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.79988513554914
$ws.Range("C2").Value = 0.164643636088581
$ws.Range("B3").Value = 0.629808504017636
$ws.Range("C3").Value = 0.206972626700285

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -3.46876419652843
$ws.Range("C2").Value = 0.299780689171852
$ws.Range("B3").Value = 0.0242905260330743
$ws.Range("C3").Value = 0.0184675218673625

# --- exp: no value changes ---

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.394501017374366
$ws.Range("B2").Value = -0.14457190032534
$ws.Range("A3").Value = -0.14457190032534
$ws.Range("B3").Value = 0.059732725393675

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.487346132809803
$ws.Range("B2").Value = -0.149652715148843
$ws.Range("A3").Value = -0.149652715148843
$ws.Range("B3").Value = 0.0492100000948854

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0271075269044691
$ws.Range("B2").Value = -0.00816904122751264
$ws.Range("A3").Value = -0.00816904122751264
$ws.Range("B3").Value = 0.0428376682032157

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0898684616003503
$ws.Range("B2").Value = -0.0041447436937204
$ws.Range("A3").Value = -0.0041447436937204
$ws.Range("B3").Value = 0.000341049363921511

# --- exp cov: no value changes ---
